# MotorMappings.xlsx edit:
#  - Rework the "1-6-19" sheet's motor list (armMotor/clawMotor/gunMotor move up one
#    row, frontSideMotor row is gone, and a Control/Button column is filled in for
#    those rows).
#  - Duplicate the (now updated) "1-6-19" sheet into a brand new "1-12-19" sheet that
#    becomes the active tab, and give that copy its own tweaked motor list (including
#    a new rubberMotor entry).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the existing "1-6-19" sheet
# ---------------------------------------------------------------------------
$sheet619 = $wb.Worksheets.Item("1-6-19")

$sheet619.Range("B4").Value = "armMotor"
$sheet619.Range("C4").Value = "5U"

$sheet619.Range("B5").Value = "clawMotor"
$sheet619.Range("C5").Value = "6U"

$sheet619.Range("B6").Value = "gunMotor"
$sheet619.Range("C6").Value = "7R"

# Row 7 no longer has a motor / control entry.
$sheet619.Range("B7").Clear()
$sheet619.Range("C7").Clear()

[void]$sheet619.Range("B7").Select()

# ---------------------------------------------------------------------------
# 2. Duplicate "1-6-19" -> new sheet "1-12-19", placed right after it, and make
#    it the active sheet/tab.
# ---------------------------------------------------------------------------
$sheet619.Copy([System.Type]::Missing, $sheet619)
$sheet1219 = $wb.Worksheets.Item($sheet619.Index + 1)
$sheet1219.Name = "1-12-19"

# ---------------------------------------------------------------------------
# 3. Update the new "1-12-19" sheet's motor list
# ---------------------------------------------------------------------------
$sheet1219.Range("B4").Value = "armMotor"
$sheet1219.Range("B5").Value = "gunMotor"
$sheet1219.Range("B6").Value = "clawMotor"
$sheet1219.Range("B7").Value = "rubberMotor"

$sheet1219.Range("C4").Value = "5U"
$sheet1219.Range("C5").Value = "6U"
$sheet1219.Range("C6").Value = "8L"
$sheet1219.Range("C7").Value = "7L/7R"

[void]$sheet1219.Range("B6").Select()

$sheet1219.Activate()
